# Refresh the cryptocurrency price / 1h-volume table with the latest scrape.
# (mirrors the "Updated cryptos list ... with GitHub Actions" workflow commit)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text updates: coin names, links, % changes, and Price strings that
#     Excel cannot mistake for a number (multiple dots, etc.) ----------------
$ws.Range("D2").Value = "43.899.05"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "2.348.39"
$ws.Range("E3").Value = "  -1.30%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  -1.31%  "
$ws.Range("E6").Value = "  -3.39%  "
$ws.Range("E7").Value = "  -6.01%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -2.73%  "
$ws.Range("E10").Value = "  -3.75%  "
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("E14").Value = "  -4.78%  "
$ws.Range("D15").Value = "2.697.90"
$ws.Range("E15").Value = "  -1.21%  "
$ws.Range("E16").Value = "  -5.09%  "
$ws.Range("E17").Value = "  -2.87%  "
$ws.Range("D18").Value = "2.349.72"
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("D19").Value = "43.788.76"
$ws.Range("E19").Value = "  -1.40%  "
$ws.Range("E20").Value = "  -1.01%  "
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("E22").Value = "  -0.70%  "
$ws.Range("E23").Value = "  -1.28%  "
$ws.Range("E24").Value = "  +7.91%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("E26").Value = "  +0.72%  "
$ws.Range("E27").Value = "  -3.38%  "
$ws.Range("E28").Value = "  -4.39%  "
$ws.Range("E29").Value = "  -3.56%  "
$ws.Range("E30").Value = "  +0.76%  "
$ws.Range("E31").Value = "  -3.36%  "
$ws.Range("E32").Value = "  -3.01%  "
$ws.Range("E33").Value = "  +0.39%  "
$ws.Range("E34").Value = "  -1.86%  "
$ws.Range("E35").Value = "  -4.58%  "
$ws.Range("E36").Value = "  +0.67%  "
$ws.Range("E37").Value = "  -3.65%  "
$ws.Range("E39").Value = "  -5.09%  "
$ws.Range("E40").Value = "  -1.11%  "
$ws.Range("E41").Value = "  +21.71%  "
$ws.Range("E42").Value = "  +17.26%  "
$ws.Range("E43").Value = "  +7.70%  "
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("E45").Value = "  -1.52%  "
$ws.Range("E46").Value = "  +1.59%  "
$ws.Range("E47").Value = "  -2.64%  "
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("E48").Value = "  -2.65%  "
$ws.Range("B49").Value = "BinanceUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").Value = "  -3.81%  "
$ws.Range("E51").Value = "  -5.72%  "

# --- Price strings that ARE valid numbers (e.g. "1.00", "240.10") -----------
#     The source sheet stores every Price cell as literal text, so plain
#     ".Value = ..." would let Excel "helpfully" reinterpret these as numbers
#     (dropping trailing zeros, changing the cell type). Force the text format
#     first, assign, then drop back to the default "Normal" style so the cell
#     keeps its original (unstyled) look with no visible formatting change.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.671"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.22"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.593"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1000"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.32"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "32.75"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.108"
$ws.Range("D13").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.32"
$ws.Range("D16").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "78.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "254.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.48"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "176.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.40"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "66.70"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "99.25"
$ws.Range("D50").Style = "Normal"
